# Apply the "Dev V2" edits to 00_nightly_only/whole_flow.xlsx (Sheet1).
#
# Real content changes (everything else in the raw OOXML diff is just the
# automatic shared-string table re-indexing caused by these edits):
#   1. A16 was blank -> "yes"
#   2. A17 was blank -> "yes"
#   3. B32 "list (of data.frames)" -> "list (of data.frames) / Table_List"
#   4. Selected cell moves from C15 to C14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "yes"
$ws.Range("A17").Value = "yes"
$ws.Range("B32").Value = "list (of data.frames) / Table_List"

$ws.Activate()
$ws.Range("C14").Select()
